$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Semester is now recognised (20222 exists) for rows 2,4,5,6,7: plain text cell, no highlight ---
$ws.Range("B2").ClearFormats()
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "20222"
$ws.Range("B2").ClearFormats()

$ws.Range("B4").ClearFormats()
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "20222"
$ws.Range("B4").ClearFormats()

$ws.Range("B5").ClearFormats()
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "20222"
$ws.Range("B5").ClearFormats()

$ws.Range("B6").ClearFormats()
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "20222"
$ws.Range("B6").ClearFormats()

$ws.Range("B7").ClearFormats()
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "20222"
$ws.Range("B7").ClearFormats()

# --- A3 row now fails the "student id not found" check -> copy highlight style from A2 ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "Mã sv không tồn tại"

# --- I4, I5, I7 rows are now flagged as duplicate data -> copy highlight style from I6 ---
$ws.Range("I6").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = "Dữ liệu của hàng này đã tồn tại"

$ws.Range("I6").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = "Dữ liệu của hàng này đã tồn tại"

$ws.Range("I6").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = "Dữ liệu của hàng này đã tồn tại"

$excel.CutCopyMode = 0

# --- Updated ĐRL (col E) / TC tích luỹ (col F) numbers recomputed by the import job ---
$ws.Range("E2").Value = 102
$ws.Range("F2").Value = 82

$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 66

$ws.Range("E4").Value = 86
$ws.Range("F4").Value = 82

$ws.Range("E5").Value = 46
$ws.Range("F5").Value = 77

$ws.Range("E6").Value = 106
$ws.Range("F6").Value = 76

$ws.Range("E7").Value = 59
$ws.Range("F7").Value = 94
